$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) hold numeric- and percent-looking
# strings. Excel auto-converts such text into real numbers on assignment,
# so those target cells are pre-formatted as Text ("@") to preserve them
# as literal strings, matching how the source data is stored.
$numericLookingRefs = @(
    "D2", "E2", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10",
    "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18",
    "E18", "D20", "E20", "D21", "E21", "D22", "E22", "E23", "E24", "D25", "E25", "E26", "E27", "D39", "E39",
    "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47"
)
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the new values cell by cell, in sheet (row-major) order.
$ws.Range("D2").Value = "305.09"
$ws.Range("E2").Value = "2.18%"
$ws.Range("E3").Value = "-0.26%"
$ws.Range("D4").Value = "5.175"
$ws.Range("E4").Value = "1.59%"
$ws.Range("D5").Value = "0.07520"
$ws.Range("E5").Value = "-0.18%"
$ws.Range("D6").Value = "2.277"
$ws.Range("E6").Value = "30.92%"
$ws.Range("D7").Value = "8.025"
$ws.Range("E7").Value = "3.35%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9163"
$ws.Range("E8").Value = "-1.25%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1740"
$ws.Range("E9").Value = "2.09%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.07661"
$ws.Range("E10").Value = "3.96%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08270"
$ws.Range("E11").Value = "3.75%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03042"
$ws.Range("E12").Value = "-0.41%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09938"
$ws.Range("E13").Value = "0.51%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001506"
$ws.Range("E14").Value = "0.75%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.006095"
$ws.Range("E15").Value = "-6.88%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.497"
$ws.Range("E16").Value = "1.37%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "3.879"
$ws.Range("E17").Value = "2.34%"
$ws.Range("D18").Value = "2.244"
$ws.Range("E18").Value = "1.15%"
$ws.Range("D20").Value = "0.1329"
$ws.Range("E20").Value = "0.97%"
$ws.Range("D21").Value = "4.650"
$ws.Range("E21").Value = "2.01%"
$ws.Range("D22").Value = "0.04609"
$ws.Range("E22").Value = "-0.81%"
$ws.Range("E23").Value = "1.04%"
$ws.Range("E24").Value = "3.78%"
$ws.Range("D25").Value = "0.004536"
$ws.Range("E25").Value = "2.60%"
$ws.Range("E26").Value = "-7.15%"
$ws.Range("E27").Value = "41.61%"
$ws.Range("D39").Value = "0.01759"
$ws.Range("E39").Value = "5.16%"
$ws.Range("D40").Value = "0.04571"
$ws.Range("E40").Value = "0.60%"
$ws.Range("D41").Value = "0.007262"
$ws.Range("E41").Value = "3.10%"
$ws.Range("D42").Value = "0.1363"
$ws.Range("E42").Value = "2.71%"
$ws.Range("D43").Value = "0.002198"
$ws.Range("E43").Value = "6.79%"
$ws.Range("D44").Value = "0.01079"
$ws.Range("E44").Value = "-15.61%"
$ws.Range("D45").Value = "0.00006533"
$ws.Range("E45").Value = "7.99%"
$ws.Range("E46").Value = "-57.22%"
$ws.Range("D47").Value = "0.009888"
$ws.Range("E47").Value = "-23.74%"
